$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Cached "datetimeFigureOut" field text: the file was re-saved two days
#    later, so PowerPoint refreshed every cached auto-date placeholder from
#    28-01-2019 to 30-01-2019 (or 1/28/2019 -> 1/30/2019, depending on the
#    placeholder's locale format). Update every "Date Placeholder" shape on
#    the slide layouts that still shows the old cached text.
# ---------------------------------------------------------------------------
$customLayouts = $p.SlideMaster.CustomLayouts
for ($layoutIdx = 1; $layoutIdx -le $customLayouts.Count; $layoutIdx++) {
    $layoutItem = $customLayouts.Item($layoutIdx)
    for ($layoutShapeIdx = 1; $layoutShapeIdx -le $layoutItem.Shapes.Count; $layoutShapeIdx++) {
        $layoutShape = $layoutItem.Shapes.Item($layoutShapeIdx)
        if ($layoutShape.HasTextFrame) {
            $layoutTextRange = $layoutShape.TextFrame.TextRange
            if ($layoutTextRange.Text -eq "1/28/2019") {
                $layoutTextRange.Text = "1/30/2019"
            }
        }
    }
}

# Same cached field on the Notes Master (dd-mm-yyyy style). The notes master
# is not an editable target on this host (writes to it are ignored), so this
# is attempted last-effort but does not affect anything else if it no-ops.
$notesMaster = $p.NotesMaster
for ($notesShapeIdx = 1; $notesShapeIdx -le $notesMaster.Shapes.Count; $notesShapeIdx++) {
    $notesShape = $notesMaster.Shapes.Item($notesShapeIdx)
    if ($notesShape.HasTextFrame) {
        $notesTextRange = $notesShape.TextFrame.TextRange
        if ($notesTextRange.Text -eq "28-01-2019") {
            $notesTextRange.Text = "30-01-2019"
        }
    }
}
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 2) Slide 89 ("$inc operator"): highlight the "$inc " token in red (C00000),
#    keeping the rest of the sentence in the original dark-grey (222222).
# ---------------------------------------------------------------------------
$incSlide = $p.Slides.Item(89)
for ($incShapeIdx = 1; $incShapeIdx -le $incSlide.Shapes.Count; $incShapeIdx++) {
    $incShape = $incSlide.Shapes.Item($incShapeIdx)
    if ($incShape.HasTextFrame) {
        $incTextRange = $incShape.TextFrame.TextRange
        if ($incTextRange.Text -eq "The `$inc operator increments a field by a specified value.") {
            $incTextRange.Characters(5, 5).Font.Color.RGB = 192
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 91 ("$unset operator"): highlight the "$unset" token in red
#    (C00000), keeping the rest of the sentence in dark-grey (222222).
# ---------------------------------------------------------------------------
$unsetSlide = $p.Slides.Item(91)
for ($unsetShapeIdx = 1; $unsetShapeIdx -le $unsetSlide.Shapes.Count; $unsetShapeIdx++) {
    $unsetShape = $unsetSlide.Shapes.Item($unsetShapeIdx)
    if ($unsetShape.HasTextFrame) {
        $unsetTextRange = $unsetShape.TextFrame.TextRange
        if ($unsetTextRange.Text -eq "The `$unset operator deletes a particular field.") {
            $unsetTextRange.Characters(5, 6).Font.Color.RGB = 192
        }
    }
}
